$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Islandora Metadata Template")

# Insert a new column before column G (the new column becomes G; the old
# G, H, I, ... shift right by one). Excel's default Insert behaviour
# copies formatting from the column to the left (F), which matches the
# target column width/style for the new column.
$ws.Columns("G").Insert()

# Make sure the new column's width matches column F exactly (Insert
# should already copy this, but set explicitly to be safe).
$ws.Columns("G").ColumnWidth = $ws.Columns("F").ColumnWidth

# New header cell for the inserted column.
$ws.Range("G1").Value = "label"

$ws.Activate()
$ws.Range("G4").Select()
